# Added email alert parser and renderer
# Insert a new "email channels" worksheet as the first sheet in the workbook,
# populate it with the email channel data, apply the header/hyperlink
# styling, and add the mailto hyperlinks - matching the authored workbook.

$wb = $excel.ActiveWorkbook

# --- Create the new sheet and put it first ---------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "email channels"

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 33.44140625
$ws.Columns.Item(2).ColumnWidth = 20.109375
$ws.Columns.Item(3).ColumnWidth = 38.88671875
$ws.Columns.Item(4).ColumnWidth = 22.88671875

# --- Header row (written first, in column order) ---------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Type"

# --- Fill the "Type" column down first (matches original authoring order) --
$ws.Range("B2").Value = "email-channel"
$ws.Range("B3").Value = "email-channel"
$ws.Range("B4").Value = "email-channel"

# --- Recipients header -------------------------------------------------------
$ws.Range("C1").Value = "Recipients"

# --- Name column data --------------------------------------------------------
$ws.Range("A2").Value = "John Doe"
$ws.Range("A3").Value = "Jane Doe"
$ws.Range("A4").Value = "Bob Smith"

# --- Recipients column data --------------------------------------------------
$ws.Range("C2").Value = "john.doe@test.com"
$ws.Range("C3").Value = "jane.doe@test.com"
$ws.Range("C4").Value = "bob.smith@test.com"

# --- Final header -------------------------------------------------------------
$ws.Range("D1").Value = "include JSON Attachment"

# --- "include JSON Attachment" boolean column --------------------------------
$ws.Range("D2").Value = $true
$ws.Range("D3").Value = $false
$ws.Range("D4").Value = $true

# --- Formatting: bold header row, wrap the Recipients header cell ----------
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("C1").WrapText = $true
$ws.Columns.Item(3).WrapText = $true

# --- Hyperlinks (mailto links) ------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:gerald_curley@hotmail.com", "", "", "gerald_curley@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:john.doe@test.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:gerald_curley@hotmail.com", "", "", "gerald_curley@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:gerald_curley@hotmail.com", "", "", "gerald_curley@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:jane.doe@test.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:bob.smith@test.com")

# --- Selection / active tab: the new sheet becomes the active tab, ----------
# selecting D2 there (matches the authored file); re-select "nrql" whole
# sheet first so that sheet keeps its own (unrelated) full-sheet selection.
$nrql = $wb.Worksheets.Item("nrql")
$null = $nrql.Cells.Select()

$null = $ws.Range("D2").Select()
